# Update automàtic: dades i banners [2026-02-14 20:50]
# Applies the latest meteo.cat XEMA station readings refresh to the daily summary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-14 20:48:28"
$ws.Range("G2").Value = "148 cm"
$ws.Range("I2").Value = "35.4 mm"
$ws.Range("E3").Value = "2026-02-14 20:48:31"
$ws.Range("I3").Value = "15.0 mm"
$ws.Range("N3").Value = "-7.3 °C 20:23 TU"
$ws.Range("O3").Value = "-5.2 °C"
$ws.Range("E4").Value = "2026-02-14 20:48:33"
$ws.Range("J4").Value = "997.1 hPa"
$ws.Range("N4").Value = "5.0 °C 20:08 TU"
$ws.Range("O4").Value = "10.7 °C"
$ws.Range("E5").Value = "2026-02-14 20:48:36"
$ws.Range("I5").Value = "22.0 mm"
$ws.Range("N5").Value = "-6.9 °C 20:11 TU"
$ws.Range("O5").Value = "-5.2 °C"
$ws.Range("E6").Value = "2026-02-14 20:48:39"
$ws.Range("H6").Value = "'75%"
$ws.Range("J6").Value = "997.2 hPa"
$ws.Range("E7").Value = "2026-02-14 20:48:41"
$ws.Range("H7").Value = "'52%"
$ws.Range("J7").Value = "997.3 hPa"
$ws.Range("E8").Value = "2026-02-14 20:48:44"
$ws.Range("H8").Value = "'62%"
$ws.Range("J8").Value = "997.2 hPa"
$ws.Range("E9").Value = "2026-02-14 20:48:46"
$ws.Range("H9").Value = "'54%"
$ws.Range("L9").Value = "63.0 km/h - 338º 20:28 TU"
$ws.Range("O9").Value = "11.8 °C"
$ws.Range("E10").Value = "2026-02-14 20:48:49"
$ws.Range("E11").Value = "2026-02-14 20:48:50"
$ws.Range("E12").Value = "2026-02-14 20:48:51"
$ws.Range("N12").Value = "10.5 °C 20:28 TU"
$ws.Range("E13").Value = "2026-02-14 20:48:52"
$ws.Range("J13").Value = "1000.0 hPa"
$ws.Range("E14").Value = "2026-02-14 20:48:53"
$ws.Range("H14").Value = "'50%"
$ws.Range("E15").Value = "2026-02-14 20:48:54"
$ws.Range("N15").Value = "9.8 °C 20:29 TU"
$ws.Range("E16").Value = "2026-02-14 20:48:55"
$ws.Range("N16").Value = "-8.9 °C 20:13 TU"
$ws.Range("O16").Value = "-6.1 °C"
$ws.Range("E17").Value = "2026-02-14 20:48:57"
$ws.Range("K17").Value = "12.4 MJ/m2"
$ws.Range("L17").Value = "60.5 km/h - 35º 20:10 TU"
$ws.Range("N17").Value = "-0.5 °C 20:27 TU"
$ws.Range("E18").Value = "2026-02-14 20:48:58"
$ws.Range("I18").Value = "1.0 mm"
$ws.Range("J18").Value = "997.4 hPa"
$ws.Range("K18").Value = "11.8 MJ/m2"
$ws.Range("L18").Value = "25.9 km/h - 306º 20:20 TU"
$ws.Range("E19").Value = "2026-02-14 20:48:59"
$ws.Range("H19").Value = "'77%"
$ws.Range("E20").Value = "2026-02-14 20:49:00"
$ws.Range("N20").Value = "-7.9 °C 20:22 TU"
$ws.Range("E21").Value = "2026-02-14 20:49:02"
$ws.Range("J21").Value = "999.9 hPa"
$ws.Range("E22").Value = "2026-02-14 20:49:05"
$ws.Range("I22").Value = "0.8 mm"
$ws.Range("N22").Value = "-9.4 °C 20:23 TU"
$ws.Range("O22").Value = "-6.9 °C"
$ws.Range("E23").Value = "2026-02-14 20:49:07"
$ws.Range("I23").Value = "39.6 mm"
$ws.Range("E24").Value = "2026-02-14 20:49:10"
$ws.Range("J24").Value = "1001.4 hPa"
$ws.Range("O24").Value = "9.5 °C"
$ws.Range("E25").Value = "2026-02-14 20:49:12"
$ws.Range("I25").Value = "15.3 mm"
$ws.Range("N25").Value = "-7.3 °C 20:17 TU"
$ws.Range("O25").Value = "-4.7 °C"
$ws.Range("E26").Value = "2026-02-14 20:49:15"
$ws.Range("E27").Value = "2026-02-14 20:49:17"
$ws.Range("O27").Value = "-3.1 °C"
$ws.Range("E28").Value = "2026-02-14 20:49:20"
$ws.Range("H28").Value = "'67%"
$ws.Range("J28").Value = "997.1 hPa"
$ws.Range("L28").Value = "55.4 km/h - 14º 20:15 TU"
$ws.Range("E29").Value = "2026-02-14 20:49:22"
$ws.Range("H29").Value = "'63%"
$ws.Range("O29").Value = "11.5 °C"
$ws.Range("E30").Value = "2026-02-14 20:49:25"
$ws.Range("J30").Value = "997.0 hPa"
$ws.Range("K30").Value = "8.1 MJ/m2"
$ws.Range("L30").Value = "107.6 km/h - 9º 20:27 TU"
$ws.Range("E31").Value = "2026-02-14 20:49:27"
$ws.Range("J31").Value = "996.1 hPa"
$ws.Range("N31").Value = "8.2 °C 20:23 TU"
$ws.Range("E32").Value = "2026-02-14 20:49:30"
$ws.Range("E33").Value = "2026-02-14 20:49:32"
$ws.Range("J33").Value = "999.3 hPa"
$ws.Range("N33").Value = "1.3 °C 20:14 TU"
$ws.Range("E34").Value = "2026-02-14 20:49:35"
$ws.Range("I34").Value = "3.5 mm"
$ws.Range("N34").Value = "-5.2 °C 20:16 TU"
$ws.Range("O34").Value = "-2.3 °C"
$ws.Range("E35").Value = "2026-02-14 20:49:38"
$ws.Range("J35").Value = "1004.0 hPa"
$ws.Range("N35").Value = "1.5 °C 20:25 TU"
$ws.Range("E36").Value = "2026-02-14 20:49:40"
$ws.Range("H36").Value = "'56%"
$ws.Range("J36").Value = "997.8 hPa"
$ws.Range("N36").Value = "10.4 °C 20:19 TU"
$ws.Range("E37").Value = "2026-02-14 20:49:43"
$ws.Range("H37").Value = "'65%"
$ws.Range("J37").Value = "998.1 hPa"
$ws.Range("E38").Value = "2026-02-14 20:49:45"
$ws.Range("L38").Value = "25.9 km/h - 329º 20:29 TU"
$ws.Range("O38").Value = "10.1 °C"
$ws.Range("E39").Value = "2026-02-14 20:49:48"
$ws.Range("I39").Value = "12.8 mm"
$ws.Range("N39").Value = "-8.5 °C 20:01 TU"
$ws.Range("O39").Value = "-5.8 °C"
$ws.Range("E40").Value = "2026-02-14 20:49:50"
$ws.Range("J40").Value = "1000.6 hPa"
$ws.Range("O40").Value = "7.0 °C"
$ws.Range("E41").Value = "2026-02-14 20:49:53"
$ws.Range("H41").Value = "'47%"
$ws.Range("J41").Value = "999.1 hPa"
$ws.Range("E42").Value = "2026-02-14 20:49:55"
$ws.Range("E43").Value = "2026-02-14 20:49:58"
$ws.Range("H43").Value = "'66%"
$ws.Range("E44").Value = "2026-02-14 20:50:00"
$ws.Range("G44").Value = "272 cm"
$ws.Range("I44").Value = "37.3 mm"
$ws.Range("N44").Value = "-7.1 °C 20:23 TU"
$ws.Range("O44").Value = "-5.4 °C"
$ws.Range("E45").Value = "2026-02-14 20:50:03"
$ws.Range("I45").Value = "13.5 mm"
$ws.Range("J45").Value = "1006.3 hPa"
$ws.Range("N45").Value = "0.7 °C 20:29 TU"
$ws.Range("E46").Value = "2026-02-14 20:50:05"
